# Token-level author diarization slides: move the two pictures on the
# "proposed approach" (Model B) slide earlier in the z-order, reflow the
# arrow/caption textbox, and tweak the wording on the "Model A" slide
# describing the trainable transformation.

$EMU_PER_POINT = 914400.0 / 72.0

function EmuToPt([double]$emu) {
    return $emu / $EMU_PER_POINT
}

# Shape.Name occasionally gets normalised by the host for auto-shapes
# (e.g. "Arrow: Right 11" reports back as "Right Arrow 11"), so look
# shapes up by their stable numeric Id instead of by name.
function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 10: "The proposed approach" (Model B)
# ---------------------------------------------------------------------
$s10 = $p.Slides.Item(10)

$picB  = Get-ShapeById $s10 3    # "Content Placeholder 2", half-size picture, rId3
$picB2 = Get-ShapeById $s10 11   # "Content Placeholder 10", quarter-size picture, rId4
$arrow = Get-ShapeById $s10 12   # "Arrow: Right 11"
$caption = Get-ShapeById $s10 13 # "TextBox 12"

# Move both pictures to the front of the shape (z-order) list, right after
# the group shape properties, ahead of the Title/body text placeholders.
$picB.ZOrder(1)    # msoSendToBack -> Content Placeholder 2 becomes shape #1
$picB2.ZOrder(3)   # msoSendBackward, repeated, to land right after picB
$picB2.ZOrder(3)
$picB2.ZOrder(3)

# Reposition the pictures slightly higher on the slide.
$picB.Top = EmuToPt 2704779
$picB2.Top = EmuToPt 2704779

# Move the right-arrow connector up along with the pictures.
$arrow.Top = EmuToPt 4257213

# Move the footnote caption textbox further down the slide.
$caption.Top = EmuToPt 6334150

# ---------------------------------------------------------------------
# Slide 9: "The proposed approach" (Model A) - wording tweaks
# ---------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$sh8 = Get-ShapeById $s9 8   # "Content Placeholder 7"
$tr = $sh8.TextFrame.TextRange

# Drop the leading "Uses " (two Croatian-tagged runs) so the remaining
# English run absorbs the whole first paragraph, then restate it.
$tr.Characters(1, 5).Delete()

$newFirst = "A trainable transformation to adapt the feature space for clustering"
$tr.Characters(1, 77).Text = $newFirst

# The second paragraph now starts right after the new first paragraph's
# text plus its paragraph break.
$secondParaStart = $newFirst.Length + 2
$tr.Characters($secondParaStart, 44).Text = "Linear (output dimension 40) or elementwise linear transformation "
